$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply changes cell-by-cell. NumberFormat is forced to Text ("@") before
# assigning values so that numeric-looking strings (e.g. "1.00", "66.972.54",
# "0.0000168") are preserved exactly as text rather than being parsed into
# floating point numbers by Excel.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '66.972.54'
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  -1.55%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.459.88'
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -2.89%  '
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '580.95'
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -1.69%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '169.17'
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -2.50%  '
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.511'
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -2.48%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '2.460.10'
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -2.91%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.134'
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -2.31%  '
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '4.88'
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -2.74%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.327'
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -4.98%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '25.16'
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -4.79%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '66.878.36'
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -1.75%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.0000168'
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -4.42%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.470.64'
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -4.71%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.90'
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -8.56%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.35'
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -8.41%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '349.58'
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -4.27%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.00'
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -3.72%  '
$c = $ws.Range('B23')
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c = $ws.Range('C23')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  -0.01%  '
$c = $ws.Range('B24')
$c.NumberFormat = "@"
$c.Value = 'Litecoin'
$c = $ws.Range('C24')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '68.62'
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -4.57%  '
$c = $ws.Range('B25')
$c.NumberFormat = "@"
$c.Value = 'NEARProtocol'
$c = $ws.Range('C25')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '4.18'
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -7.66%  '
$c = $ws.Range('B26')
$c.NumberFormat = "@"
$c.Value = 'SuiNetwork'
$c = $ws.Range('C26')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '1.79'
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -4.83%  '
$c = $ws.Range('B27')
$c.NumberFormat = "@"
$c.Value = 'Aptos'
$c = $ws.Range('C27')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.08'
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -7.89%  '
$c = $ws.Range('B28')
$c.NumberFormat = "@"
$c.Value = 'Binance-PegBSC-USD'
$c = $ws.Range('C28')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -22.89%  '
$c = $ws.Range('B29')
$c.NumberFormat = "@"
$c.Value = 'WrappedeETH'
$c = $ws.Range('C29')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.583.25'
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -3.21%  '
$c = $ws.Range('B30')
$c.NumberFormat = "@"
$c.Value = 'Bittensor'
$c = $ws.Range('C30')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '509.91'
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -4.10%  '
$c = $ws.Range('B31')
$c.NumberFormat = "@"
$c.Value = 'PEPE'
$c = $ws.Range('C31')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.0₃0891'
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -6.75%  '
$c = $ws.Range('B32')
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c = $ws.Range('C32')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.58'
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -8.26%  '
$c = $ws.Range('B33')
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c = $ws.Range('C33')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.75'
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -5.96%  '
$c = $ws.Range('B34')
$c.NumberFormat = "@"
$c.Value = 'Fetch.AI'
$c = $ws.Range('C34')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.22'
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -6.20%  '
$c = $ws.Range('B35')
$c.NumberFormat = "@"
$c.Value = 'FirstDigitalUSD'
$c = $ws.Range('C35')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c = $ws.Range('B36')
$c.NumberFormat = "@"
$c.Value = 'Monero'
$c = $ws.Range('C36')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '158.25'
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -0.99%  '
$c = $ws.Range('B37')
$c.NumberFormat = "@"
$c.Value = 'Kaspa'
$c = $ws.Range('C37')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.113'
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -11.75%  '
$c = $ws.Range('B38')
$c.NumberFormat = "@"
$c.Value = 'WhiteBITCoin'
$c = $ws.Range('C38')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '18.64'
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c = $ws.Range('B39')
$c.NumberFormat = "@"
$c.Value = 'EthereumClassic'
$c = $ws.Range('C39')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '18.15'
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -5.88%  '
$c = $ws.Range('B40')
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c = $ws.Range('C40')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.32'
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -8.05%  '
$c = $ws.Range('B41')
$c.NumberFormat = "@"
$c.Value = 'USDe'
$c = $ws.Range('C41')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c = $ws.Range('B42')
$c.NumberFormat = "@"
$c.Value = 'Stacks'
$c = $ws.Range('C42')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.66'
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -6.76%  '
$c = $ws.Range('B43')
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c = $ws.Range('C43')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '4.74'
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -7.09%  '
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -6.96%  '
$c = $ws.Range('B45')
$c.NumberFormat = "@"
$c.Value = 'dogwifhat'
$c = $ws.Range('C45')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.35'
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -4.99%  '
$c = $ws.Range('B46')
$c.NumberFormat = "@"
$c.Value = 'OKB'
$c = $ws.Range('C46')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '38.55'
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -2.13%  '
$c = $ws.Range('B47')
$c.NumberFormat = "@"
$c.Value = 'Aave'
$c = $ws.Range('C47')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '140.84'
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -5.10%  '
$c = $ws.Range('B48')
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c = $ws.Range('C48')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '3.42'
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -7.63%  '
$c = $ws.Range('B49')
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c = $ws.Range('C49')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.506'
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -8.37%  '
$c = $ws.Range('B50')
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c = $ws.Range('C50')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0727'
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -2.57%  '
$c = $ws.Range('B51')
$c.NumberFormat = "@"
$c.Value = 'BabyDogeCoin'
$c = $ws.Range('C51')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0₆0247'
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -9.97%  '
